# Make demo data more realistic and complete:
# - header row stays tag_id / variable_id
# - replace the two sample data rows with a single, more realistic row
# - drop the now-unused third row entirely

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "variable_id"

$ws.Range("A2").Value = "localisation"
$ws.Range("B2").Value = "ser_pub_loc___canton"

# Clear what used to be the third data row so the sheet shrinks to A1:B2
$ws.Range("A3").Value = ""
$ws.Range("B3").Value = ""

# Restore the bottom-right pane's active selection to B5 (matches the
# author's recorded view state after trimming the data)
$ws.Range("B5").Select() | Out-Null
